$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.743740081787109
$ws.Range("B1").Value = 3.725407838821411
$ws.Range("C1").Value = 2.02006983757019
$ws.Range("D1").Value = 1.411194801330566
$ws.Range("E1").Value = 1.204880237579346
